# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45212 (2023-10-13) to 45221 (2023-10-22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45212) {
        $cell.Value2 = 45221
    }
}
